$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the Price column as Text first so numeric-looking strings
# (e.g. "1.002") are stored as text, matching the original inline-string cells,
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "23.121.44"
$ws.Range("E2").Value = "  -3.59%  "
$ws.Range("D3").Value = "1.602.53"
$ws.Range("E3").Value = "  -2.94%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "1.001"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("D6").Value = "301.33"
$ws.Range("E6").Value = "  -2.82%  "
$ws.Range("D7").Value = "0.3779"
$ws.Range("E7").Value = "  -3.03%  "
$ws.Range("D8").Value = "0.3646"
$ws.Range("E8").Value = "  -4.58%  "
$ws.Range("D9").Value = "49.96"
$ws.Range("E9").Value = "  -4.66%  "
$ws.Range("D10").Value = "1.260"
$ws.Range("E10").Value = "  -6.50%  "
$ws.Range("B11").Value = "BinanceUSD"
$ws.Range("C11").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D11").Value = "1.002"
$ws.Range("E11").Value = "  +0.00%  "
$ws.Range("B12").Value = "Dogecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D12").Value = "0.08147"
$ws.Range("E12").Value = "  -3.46%  "
$ws.Range("E13").Value = "  -3.21%  "
$ws.Range("D14").Value = "6.630"
$ws.Range("E14").Value = "  -6.19%  "
$ws.Range("D15").Value = "7.411"
$ws.Range("E15").Value = "  -7.40%  "
$ws.Range("E16").Value = "  -4.10%  "
$ws.Range("D17").Value = "1.603.87"
$ws.Range("E17").Value = "  -2.74%  "
$ws.Range("D18").Value = "91.71"
$ws.Range("E18").Value = "  -3.04%  "
$ws.Range("D19").Value = "0.06851"
$ws.Range("E19").Value = "  -2.12%  "
$ws.Range("D20").Value = "18.25"
$ws.Range("E20").Value = "  -7.23%  "
$ws.Range("D21").Value = "6.581"
$ws.Range("E21").Value = "  -5.64%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").Value = "12.99"
$ws.Range("E23").Value = "  -5.88%  "
$ws.Range("D24").Value = "23.136.66"
$ws.Range("E25").Value = "  -4.02%  "
$ws.Range("D26").Value = "2.792"
$ws.Range("E26").Value = "  -5.97%  "
$ws.Range("D27").Value = "21.08"
$ws.Range("E27").Value = "  -4.53%  "
$ws.Range("D28").Value = "150.37"
$ws.Range("E28").Value = "  -1.30%  "
$ws.Range("D29").Value = "5.285"
$ws.Range("E29").Value = "  -2.46%  "
$ws.Range("D30").Value = "131.90"
$ws.Range("E30").Value = "  -4.51%  "
$ws.Range("D31").Value = "2.418"
$ws.Range("E31").Value = "  -3.88%  "
$ws.Range("D32").Value = "6.864"
$ws.Range("E32").Value = "  -13.55%  "
$ws.Range("D33").Value = "1.779.28"
$ws.Range("E33").Value = "  -2.68%  "
$ws.Range("D34").Value = "0.07692"
$ws.Range("E34").Value = "  -4.88%  "
$ws.Range("D35").Value = "0.9460"
$ws.Range("E35").Value = "  -7.58%  "
$ws.Range("D36").Value = "0.02768"
$ws.Range("E36").Value = "  -5.60%  "
$ws.Range("D37").Value = "6.259"
$ws.Range("E37").Value = "  -6.98%  "
$ws.Range("D38").Value = "0.2543"
$ws.Range("D39").Value = "0.08912"
$ws.Range("D40").Value = "10.11"
$ws.Range("E40").Value = "  -5.39%  "
$ws.Range("E41").Value = "  -1.89%  "
$ws.Range("D42").Value = "12.78"
$ws.Range("E42").Value = "  -4.50%  "
$ws.Range("D43").Value = "0.7105"
$ws.Range("E43").Value = "  -6.46%  "
$ws.Range("D44").Value = "15.46"
$ws.Range("E44").Value = "  -5.51%  "
$ws.Range("D45").Value = "0.6627"
$ws.Range("E45").Value = "  -4.83%  "
$ws.Range("D46").Value = "1.000"
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").Value = "2.302"
$ws.Range("E47").Value = "  -6.52%  "
$ws.Range("E48").Value = "  -3.10%  "
$ws.Range("D49").Value = "132.53"
$ws.Range("E49").Value = "  -1.67%  "
$ws.Range("D50").Value = "0.07952"
$ws.Range("E50").Value = "  -4.57%  "
$ws.Range("D51").Value = "1.211"
$ws.Range("E51").Value = "  -0.90%  "

# Restore the column's number format back to General so the cell style
# matches the original (unstyled) cells exactly.
$ws.Range("D2:D51").NumberFormat = "General"
$ws.Range("D2:D51").Style = "Normal"
